# Zeitplanung.xlsx update - "Wahrscheinlich letzte Aktualisierung des Zeitplans (vh)"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgabenliste Projekt 1")

# Row 38: "tatsächliche Fertigstellung" (H) now filled in with a date
$ws.Range("H38").Value = [DateTime]"2015-11-28"

# Row 40: "Fortschritt" (F) now marked 1 (100%), matching E40
$ws.Range("F40").Value = 1

# Row 50: task is now fully done -> % erledigt (E) and Fortschritt (F) = 1, and
# "tatsächliche Fertigstellung" (H) date filled in
$ws.Range("E50").Value = 1
$ws.Range("F50").Value = 1
$ws.Range("H50").Value = [DateTime]"2015-12-12"

# Row 51: % erledigt (E) cleared (no longer 0%), and "tatsächliche Fertigstellung" (H)
# marked with "/" (not applicable / cancelled)
$ws.Range("E51").ClearContents()
$ws.Range("H51").Value = "/"

# Row 52: "tatsächliche Fertigstellung" (H) now filled in with a date
$ws.Range("H52").Value = [DateTime]"2015-12-13"

# Update view state to match where the author was last working
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("F51").Select()
